$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice2")

# Update the surviving data row (row 6) with the new values.
$ws.Range("H6").Value = "SO-0017943"
$ws.Range("L6").Value = 25000

# Remove the three now-obsolete data rows (old rows 7-9); this shifts the
# Total row (old row 10) up to row 7 and everything below it up as well.
$ws.Rows("7:9").Delete()

# Leave the selection where the editor ended up.
[void]$ws.Range("M14").Select()
